# Swap the order of "System" and the email address in the
# "Recorded By" column (G) wherever both recorders are listed as
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldVal = "System, dnasr281@gmail.com"
$newVal = "dnasr281@gmail.com, System"

$searchRange = $ws.Range("G1:G319")

$firstAddress = $null
$found = $searchRange.Find($oldVal)

while ($found -ne $null) {
    if ($firstAddress -eq $null) {
        $firstAddress = $found.Address()
    } elseif ($found.Address() -eq $firstAddress) {
        break
    }

    $found.Value2 = $newVal

    $found = $searchRange.FindNext($found)
}
